# Apply cell value updates from the crypto price refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.980.52"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.555.80"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.97"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.777.96"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "1.555.68"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "26.989.52"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.80"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.20"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.33"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.97"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "1.423.66"
$ws.Range("E33").Value = "  +5.05%  "
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("E35").Value = "  +3.98%  "
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.812"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +4.90%  "
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.56"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "1.691.40"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "0.0₇0995"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  +1.17%  "
